# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns.
# Cells that would otherwise auto-parse as numbers are forced back to text
# with a leading apostrophe, matching the source data's text-typed cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.063.41"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.639.39"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'213.94"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "'0.5243"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("D8").Value = "'0.2599"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "'0.06296"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'20.61"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").Value = "'0.07679"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "1.649.03"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "'4.400"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").Value = "1.857.88"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").Value = "'0.5518"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "0.0₅8187"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "'64.92"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "26.052.86"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D20").Value = "'4.682"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "'188.43"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").Value = "'10.18"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").Value = "'6.158"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "'145.31"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").Value = "'0.1206"
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").Value = "'7.396"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "'1.381"
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("E30").Value = "  -5.95%  "
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").Value = "'3.427"
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").Value = "'3.394"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "'1.645"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "'0.9816"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").Value = "'2.393"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "'2.759"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").Value = "'0.5620"
$ws.Range("E38").Value = "  -6.21%  "
$ws.Range("D39").Value = "'0.01614"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'0.8470"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D42").Value = "'5.686"
$ws.Range("E42").Value = "  -6.78%  "
$ws.Range("D43").Value = "1.024.57"
$ws.Range("E43").Value = "  -7.36%  "
$ws.Range("D44").Value = "'100.11"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "1.786.81"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").Value = "0.0₈106"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").Value = "'55.72"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").Value = "'0.9998"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "'8.013"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "'0.4215"
